$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "10+1="
$t.Cell(1,2).Range.Text = "22+29="
$t.Cell(1,3).Range.Text = "69-38="
$t.Cell(1,4).Range.Text = "4+78="
$t.Cell(1,5).Range.Text = "96-84="
$t.Cell(2,1).Range.Text = "54-32="
$t.Cell(2,2).Range.Text = "50-19="
$t.Cell(2,3).Range.Text = "87-74="
$t.Cell(2,4).Range.Text = "72+10="
$t.Cell(2,5).Range.Text = "37-17="
$t.Cell(3,1).Range.Text = "65+24="
$t.Cell(3,2).Range.Text = "67-41="
$t.Cell(3,3).Range.Text = "35-7="
$t.Cell(3,4).Range.Text = "98-6="
$t.Cell(3,5).Range.Text = "45-6="
$t.Cell(4,1).Range.Text = "97-73="
$t.Cell(4,2).Range.Text = "33+62="
$t.Cell(4,3).Range.Text = "48-25="
$t.Cell(4,4).Range.Text = "43-19="
$t.Cell(4,5).Range.Text = "19+78="
$t.Cell(5,1).Range.Text = "41-7="
$t.Cell(5,2).Range.Text = "93-16="
$t.Cell(5,3).Range.Text = "55+19="
$t.Cell(5,4).Range.Text = "28+45="
$t.Cell(5,5).Range.Text = "23+44="
$t.Cell(6,1).Range.Text = "86-58="
$t.Cell(6,2).Range.Text = "44-21="
$t.Cell(6,3).Range.Text = "42-36="
$t.Cell(6,4).Range.Text = "10-0="
$t.Cell(6,5).Range.Text = "39+15="
$t.Cell(7,1).Range.Text = "95-71="
$t.Cell(7,2).Range.Text = "87-66="
$t.Cell(7,3).Range.Text = "53-30="
$t.Cell(7,4).Range.Text = "35-25="
$t.Cell(7,5).Range.Text = "57+33="
$t.Cell(8,1).Range.Text = "67-10="
$t.Cell(8,2).Range.Text = "44+1="
$t.Cell(8,3).Range.Text = "27+41="
$t.Cell(8,4).Range.Text = "89-78="
$t.Cell(8,5).Range.Text = "21+67="
$t.Cell(9,1).Range.Text = "23+19="
$t.Cell(9,2).Range.Text = "75-26="
$t.Cell(9,3).Range.Text = "10+13="
$t.Cell(9,4).Range.Text = "15+17="
$t.Cell(9,5).Range.Text = "34-29="
$t.Cell(10,1).Range.Text = "78-66="
$t.Cell(10,2).Range.Text = "36+15="
$t.Cell(10,3).Range.Text = "12+27="
$t.Cell(10,4).Range.Text = "76-9="
$t.Cell(10,5).Range.Text = "37-4="
$t.Cell(11,1).Range.Text = "44-30="
$t.Cell(11,2).Range.Text = "50+32="
$t.Cell(11,3).Range.Text = "83-60="
$t.Cell(11,4).Range.Text = "87-63="
$t.Cell(11,5).Range.Text = "39-18="
$t.Cell(12,1).Range.Text = "52-12="
$t.Cell(12,2).Range.Text = "81-77="
$t.Cell(12,3).Range.Text = "23-15="
$t.Cell(12,4).Range.Text = "67-35="
$t.Cell(12,5).Range.Text = "38-19="
$t.Cell(13,1).Range.Text = "70-55="
$t.Cell(13,2).Range.Text = "37-32="
$t.Cell(13,3).Range.Text = "29+58="
$t.Cell(13,4).Range.Text = "94-24="
$t.Cell(13,5).Range.Text = "65+16="
$t.Cell(14,1).Range.Text = "85-16="
$t.Cell(14,2).Range.Text = "64+26="
$t.Cell(14,3).Range.Text = "75-38="
$t.Cell(14,4).Range.Text = "37+49="
$t.Cell(14,5).Range.Text = "75+15="
$t.Cell(15,1).Range.Text = "53-32="
$t.Cell(15,2).Range.Text = "15+76="
$t.Cell(15,3).Range.Text = "47+15="
$t.Cell(15,4).Range.Text = "36-25="
$t.Cell(15,5).Range.Text = "82-35="
$t.Cell(16,1).Range.Text = "37+38="
$t.Cell(16,2).Range.Text = "28+45="
$t.Cell(16,3).Range.Text = "10+23="
$t.Cell(16,4).Range.Text = "73-6="
$t.Cell(16,5).Range.Text = "97-23="
$t.Cell(17,1).Range.Text = "91-4="
$t.Cell(17,2).Range.Text = "64-4="
$t.Cell(17,3).Range.Text = "58-3="
$t.Cell(17,4).Range.Text = "88-23="
$t.Cell(17,5).Range.Text = "24+35="
$t.Cell(18,1).Range.Text = "95-44="
$t.Cell(18,2).Range.Text = "47+51="
$t.Cell(18,3).Range.Text = "46-38="
$t.Cell(18,4).Range.Text = "36+55="
$t.Cell(18,5).Range.Text = "62-11="
$t.Cell(19,1).Range.Text = "4+13="
$t.Cell(19,2).Range.Text = "10+35="
$t.Cell(19,3).Range.Text = "33+8="
$t.Cell(19,4).Range.Text = "87-12="
$t.Cell(19,5).Range.Text = "34-22="
$t.Cell(20,1).Range.Text = "97-71="
$t.Cell(20,2).Range.Text = "58+16="
$t.Cell(20,3).Range.Text = "93+0="
$t.Cell(20,4).Range.Text = "94-73="
$t.Cell(20,5).Range.Text = "39-25="
